$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "O.Z TOPICAL LOTION 120ML" (row 73)
$ws.Rows(73).Delete()

# Update the grand total (now at row 131 after the delete) to reflect the removed item's price
$ws.Range("P131").Value = 9342.9449999999997

# Update the generated-on timestamp text (now at row 132 after the delete, merged A132:F132)
$ws.Range("A132").Value = "Monday, 6 October, 2025 8:18 PM"
